# Auto-generated Excel COM-interop script applying cached-value updates
# described in the commit diff for Spriggan_Profits.xlsx (one row of updated
# market-board derived values per hunk, across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 125
$ws.Range("I5").Value = 125
$ws.Range("K5").Value = 125
$ws.Range("M5").Value = -10

# Row 33
$ws.Range("H33").Value = 74379.47
$ws.Range("J33").Value = 492.42856
$ws.Range("L33").Value = 492.42856
$ws.Range("N33").Value = -950.4285600000001

# Row 40
$ws.Range("H40").Value = 5933.3335
$ws.Range("I40").Value = 4333.3335
$ws.Range("J40").Value = 7533.3335
$ws.Range("K40").Value = 4333.3335
$ws.Range("L40").Value = 7533.3335
$ws.Range("M40").Value = -4158.3335
$ws.Range("N40").Value = -7883.3335

# Row 51
$ws.Range("H51").Value = 70000
$ws.Range("J51").Value = 70000
$ws.Range("L51").Value = 70000
$ws.Range("N51").Value = -70968

# Row 58
$ws.Range("H58").Value = 4091.3333
$ws.Range("I58").Value = 762.25
$ws.Range("J58").Value = 10749.5
$ws.Range("K58").Value = 2286.75
$ws.Range("L58").Value = 32248.5
$ws.Range("M58").Value = -2136.75
$ws.Range("N58").Value = -32548.5

# Row 70
$ws.Range("H70").Value = 4558
$ws.Range("I70").Value = 2731.2727
$ws.Range("J70").Value = 6232.5
$ws.Range("K70").Value = 8193.8181
$ws.Range("L70").Value = 18697.5
$ws.Range("M70").Value = -7923.8181
$ws.Range("N70").Value = -19237.5

# Row 73
$ws.Range("H73").Value = 4558
$ws.Range("I73").Value = 2731.2727
$ws.Range("J73").Value = 6232.5
$ws.Range("K73").Value = 8193.8181
$ws.Range("L73").Value = 18697.5
$ws.Range("M73").Value = -7257.8181
$ws.Range("N73").Value = -20569.5

# Row 141
$ws.Range("H141").Value = 3211
$ws.Range("I141").Value = 3312.5715
$ws.Range("K141").Value = 9937.7145
$ws.Range("M141").Value = -4757.7145

$ws = $wb.Sheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 7694246
$ws.Range("I102").Value = 8335183
$ws.Range("K102").Value = 8335183
$ws.Range("M102").Value = -8333561

$ws = $wb.Sheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 2150.2942
$ws.Range("I22").Value = 2587.2
$ws.Range("J22").Value = 1526.1428
$ws.Range("K22").Value = 2587.2
$ws.Range("L22").Value = 1526.1428
$ws.Range("M22").Value = -2414.2
$ws.Range("N22").Value = -1872.1428

# Row 86
$ws.Range("H86").Value = 3781.1667
$ws.Range("I86").Value = 3942.125
$ws.Range("J86").Value = 3652.4
$ws.Range("K86").Value = 3942.125
$ws.Range("L86").Value = 3652.4
$ws.Range("M86").Value = -2819.125
$ws.Range("N86").Value = -5898.4

# Row 89
$ws.Range("H89").Value = 3781.1667
$ws.Range("I89").Value = 3942.125
$ws.Range("J89").Value = 3652.4
$ws.Range("K89").Value = 19710.625
$ws.Range("L89").Value = 18262
$ws.Range("M89").Value = -14094.625
$ws.Range("N89").Value = -29494

# Row 105
$ws.Range("H105").Value = 1963.2727
$ws.Range("I105").Value = 2031.125
$ws.Range("J105").Value = 1782.3334
$ws.Range("K105").Value = 2031.125
$ws.Range("L105").Value = 1782.3334
$ws.Range("M105").Value = -284.125
$ws.Range("N105").Value = -5276.3334

$ws = $wb.Sheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 6641.5293
$ws.Range("I22").Value = 7861.857
$ws.Range("J22").Value = 946.6667
$ws.Range("K22").Value = 7861.857
$ws.Range("L22").Value = 946.6667
$ws.Range("M22").Value = -7511.857
$ws.Range("N22").Value = -1646.6667

# Row 54
$ws.Range("H54").Value = 39999.5
$ws.Range("J54").Value = 39999.5
$ws.Range("L54").Value = 39999.5
$ws.Range("N54").Value = -41315.5

# Row 105
$ws.Range("H105").Value = 2552328
$ws.Range("I105").Value = 3402579
$ws.Range("K105").Value = 3402579
$ws.Range("M105").Value = -3400832

# Row 132
$ws.Range("H132").Value = 83337030
$ws.Range("I132").Value = 100003640
$ws.Range("K132").Value = 300010920
$ws.Range("M132").Value = -300008390

$ws = $wb.Sheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354

$ws = $wb.Sheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 159.5
$ws.Range("J2").Value = 223.16667
$ws.Range("L2").Value = 223.16667
$ws.Range("N2").Value = -449.16667

# Row 38
$ws.Range("H38").Value = 29999
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 29999
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 29999
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -30925

# Row 40
$ws.Range("H40").Value = 29900
$ws.Range("I40").Value = 29900
$ws.Range("K40").Value = 29900
$ws.Range("M40").Value = -29749

# Row 46
$ws.Range("H46").Value = 5899.8
$ws.Range("I46").Value = 3166.6667
$ws.Range("K46").Value = 3166.6667
$ws.Range("M46").Value = -3010.6667

# Row 57
$ws.Range("H57").Value = 19950
$ws.Range("I57").Value = 19950
$ws.Range("K57").Value = 19950
$ws.Range("M57").Value = -19130

# Row 111
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -56134

# Row 122
$ws.Range("H122").Value = 4554.4116
$ws.Range("J122").Value = 5749.3335
$ws.Range("L122").Value = 17248.0005
$ws.Range("N122").Value = -22148.0005

$ws = $wb.Sheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2328
$ws.Range("I22").Value = 2101.875
$ws.Range("K22").Value = 2101.875
$ws.Range("M22").Value = -1806.875

# Row 27
$ws.Range("H27").Value = 2328
$ws.Range("I27").Value = 2101.875
$ws.Range("K27").Value = 2101.875
$ws.Range("M27").Value = -1994.875

# Row 40
$ws.Range("H40").Value = 3899.2666
$ws.Range("I40").Value = 3966.3572
$ws.Range("J40").Value = 2960
$ws.Range("K40").Value = 3966.3572
$ws.Range("L40").Value = 2960
$ws.Range("M40").Value = -3830.3572
$ws.Range("N40").Value = -3232

$ws = $wb.Sheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1294.875
$ws.Range("I96").Value = 1276.5
$ws.Range("J96").Value = 1350
$ws.Range("K96").Value = 1276.5
$ws.Range("L96").Value = 1350
$ws.Range("M96").Value = 96.5
$ws.Range("N96").Value = -4096

# Row 100
$ws.Range("H100").Value = 3105.8
$ws.Range("I100").Value = 3195.2307
$ws.Range("J100").Value = 2524.5
$ws.Range("K100").Value = 6390.4614
$ws.Range("L100").Value = 5049
$ws.Range("M100").Value = -5849.4614
$ws.Range("N100").Value = -6131

# Row 136
$ws.Range("H136").Value = 29414468
$ws.Range("I136").Value = 33336132
$ws.Range("J136").Value = 1989.5
$ws.Range("K136").Value = 100008396
$ws.Range("L136").Value = 5968.5
$ws.Range("M136").Value = -100005846
$ws.Range("N136").Value = -11068.5
